# Update Agt-Agtr2 sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5, 6 and 7 — the table shrinks from 6 data rows to 3.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

# Row 2: FAPs -> Agt -> Agtr2 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.614029333333333
$ws.Range("H2").Value = 10.842088
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01476566666666667
$ws.Range("N2").Value = 0.044297
$ws.Range("O2").Value = 0.003369932031170864
$ws.Range("P2").Value = 0.003537680284170835
$ws.Range("Q2").Value = 0.05336355245955556
$ws.Range("R2").Value = 0.480271972136
$ws.Range("S2").Value = 0.003369932031170864
$ws.Range("T2").Value = 0.003537680284170835

# Row 3: FAPs -> Agt -> Agtr2 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.614029333333333
$ws.Range("H3").Value = 10.842088
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.743532666666667
$ws.Range("N3").Value = 11.230598
$ws.Range("O3").Value = 0.8543773151546027
$ws.Range("P3").Value = 0.8969064524470826
$ws.Range("Q3").Value = 13.52923686762489
$ws.Range("R3").Value = 121.763131808624
$ws.Range("S3").Value = 0.8543773151546027
$ws.Range("T3").Value = 0.8969064524470826

# Row 4: FAPs -> Agt -> Agtr2 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.614029333333333
$ws.Range("H4").Value = 10.842088
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6232934999999999
$ws.Range("N4").Value = 1.246587
$ws.Range("O4").Value = 0.1422527528142264
$ws.Range("P4").Value = 0.09955586726874661
$ws.Range("Q4").Value = 2.252600992276
$ws.Range("R4").Value = 13.515605953656
$ws.Range("S4").Value = 0.1422527528142264
$ws.Range("T4").Value = 0.09955586726874661
